$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A5").Value = "Further education and skills learner achievements and participation by age group"
$ws.Range("A6").Value = "Further education and skills sim achievements by age, sector subject area and level"
$ws.Range("B6").Value = "<a href='https://explore-education-statistics.service.gov.uk/find-statistics/further-education-and-skills'>Individualised Learner Record</a>"

$ws.Range("B6").Select()
